# edit.ps1 — apply "edit player loop bug fix" commit to the match scorecard sheet.
#
# The underlying bug: batting/bowling rows in the source data were associated
# with the wrong player (an off-by-one style loop bug), so batting stats
# (runs/balls/mode-of-dismissal/bowler) and bowling figures were shifted
# against the wrong names. This script rewrites the affected cells with the
# corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: David Warner
$ws.Range("B2").Value = 7
$ws.Range("C2").Value = 6
$ws.Range("D2").Value = 'Caught'
$ws.Range("E2").Value = ' Maheesh Theekshana'
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 1
$ws.Range("N2").Value = ' Josh Hazlewood'

# Row 3: Aaron Finch(C)
$ws.Range("B3").Value = 67
$ws.Range("C3").Value = 30
$ws.Range("D3").Value = 'Caught'
$ws.Range("E3").Value = ' Nuwan Pradeep'
$ws.Range("K3").Value = 16
$ws.Range("L3").Value = 8
$ws.Range("M3").Value = 'LBW'
$ws.Range("N3").Value = ' Adam Zampa'

# Row 4: Mitchell Marsh
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 'Bowled'
$ws.Range("E4").Value = ' Maheesh Theekshana'
$ws.Range("K4").Value = 12
$ws.Range("L4").Value = 5
$ws.Range("M4").Value = 'Caught'
$ws.Range("N4").Value = ' Pat Cummins'

# Row 5: Steve Smith
$ws.Range("B5").Value = 40
$ws.Range("C5").Value = 12
$ws.Range("D5").Value = 'Caught'
$ws.Range("K5").Value = 7
$ws.Range("L5").Value = 3
$ws.Range("M5").Value = 'LBW'
$ws.Range("N5").Value = ' Mitchell Starc'

# Row 6: Glenn Maxwell
$ws.Range("B6").Value = 49
$ws.Range("C6").Value = 12
$ws.Range("D6").Value = 'Bowled'
$ws.Range("E6").Value = ' Chamika Karunarathne'
$ws.Range("K6").Value = 8
$ws.Range("L6").Value = 8
$ws.Range("N6").Value = ' Josh Hazlewood'

# Row 7: Matthew Wade
$ws.Range("C7").Value = 4
$ws.Range("D7").Value = 'NOT OUT'
$ws.Range("E7").Value = ' '
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 'LBW'

# Row 8: Marcus Stionis
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = 'LBW'
$ws.Range("E8").Value = ' Maheesh Theekshana'
$ws.Range("K8").Value = 7
$ws.Range("L8").Value = 4
$ws.Range("M8").Value = 'Caught'
$ws.Range("N8").Value = ' Marcus Stionis'

# Row 9: Pat Cummins
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 'Caught'
$ws.Range("E9").Value = ' Maheesh Theekshana'
$ws.Range("K9").Value = 25
$ws.Range("L9").Value = 9
$ws.Range("M9").Value = 'Bowled'
$ws.Range("N9").Value = ' Mitchell Starc'

# Row 10: Mitchell Starc
$ws.Range("B10").Value = 17
$ws.Range("C10").Value = 6
$ws.Range("E10").Value = ' Dushmantha Chameera'
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 'Caught'
$ws.Range("N10").Value = ' Adam Zampa'

# Row 11: Adam Zampa
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 'Caught'
$ws.Range("E11").Value = ' Dushmantha Chameera'
$ws.Range("K11").Value = 28
$ws.Range("L11").Value = 12
$ws.Range("M11").Value = 'NOT OUT'
$ws.Range("N11").Value = ' '

# Row 12: Josh Hazlewood
$ws.Range("B12").Value = 11
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 'LBW'
$ws.Range("E12").Value = ' Chamika Karunarathne'
$ws.Range("K12").Value = 21
$ws.Range("L12").Value = 8
$ws.Range("M12").Value = 'LBW'

# Row 16: Innings totals
$ws.Range("A16").Value = 208
$ws.Range("C16").Value = '''13.2'
$ws.Range("D16").Value = 80
$ws.Range("J16").Value = 124
$ws.Range("L16").Value = '''10.0'
$ws.Range("M16").Value = 60

# Row 21: Bowling row 21
$ws.Range("A21").Value = 'Wanindu Hasaranga'
$ws.Range("C21").Value = 36
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 18
$ws.Range("L21").Value = 27
$ws.Range("N21").Value = 13.5

# Row 22: Bowling row 22
$ws.Range("A22").Value = 'Nuwan Pradeep'
$ws.Range("B22").Value = '''3.0'
$ws.Range("C22").Value = 39
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 13
$ws.Range("L22").Value = 23
$ws.Range("M22").Value = 2
$ws.Range("N22").Value = 11.5

# Row 23: Bowling row 23
$ws.Range("A23").Value = 'Maheesh Theekshana'
$ws.Range("B23").Value = '''3.0'
$ws.Range("C23").Value = 52
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = 17.33
$ws.Range("L23").Value = 26
$ws.Range("M23").Value = 2
$ws.Range("N23").Value = 13

# Row 24: Bowling row 24
$ws.Range("A24").Value = 'Dushmantha Chameera'
$ws.Range("C24").Value = 43
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = 14.33
$ws.Range("L24").Value = 21
$ws.Range("N24").Value = 10.5

# Row 25: Bowling row 25
$ws.Range("A25").Value = 'Chamika Karunarathne'
$ws.Range("B25").Value = '''2.2'
$ws.Range("C25").Value = 38
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 17.27
$ws.Range("K25").Value = '''2.0'
$ws.Range("L25").Value = 27
$ws.Range("M25").Value = 2
$ws.Range("N25").Value = 13.5
